$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update StructureDefinition publishing metadata ---
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date bump
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now populated
$meta.Range("B9").Value = "Alvearie Team"

# Remove the duplicated "Contact" / "No display for ContactDetail" row (row 11),
# then turn the remaining row 10 into "Jurisdiction" / "United States of America"
$meta.Rows.Item(11).Delete()
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Sheet "Elements": root Extension row now shows the real Short/Definition ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Communication Priority Sequence"
$elements.Range("L2").Value = "Priority of the communication request (1 = highest)"
